$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# RuleTable "Import" row: drop the Application pojo, keep only the Insured pojo
$ws.Range("B2").Value = "com.redhat.prudential_poc.pojo.Insured"

# Condition field header: "id" -> "insuredId"
$ws.Range("B8").Value = "insuredId"
$ws.Range("B8").Font.Name = "Arial"

# Move the active selection to B2 (matches the saved sheet view state)
$ws.Range("B2").Select()
